$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "variavel"
$ws.Range("A1").Copy()
$ws.Range("B1").PasteSpecial(-4122)
$ws.Range("B1").Value = "impacto"

$labels = @(
    "sem instrução ou fundamental incompleto",
    "proporção de pessoas com acesso simultâneo aos três serviços de saneamento básico(%)",
    "ensino fundamental completo ou médio incompleto",
    "rendimento-hora médio real habitual do trabalho principal(r$/hora)",
    "rendimento médio real habitual de todos os trabalhos(r$/mês)",
    "rendimento médio real habitual do trabalho principal(r$/mês)",
    "número de beneficiários de plano de saúde",
    "taxa de analfabetismo",
    "ensino superior completo",
    "domicílio cedido por empregador",
    "ensino médio completo ou superior incompleto",
    "outra forma",
    "taxa de formalização",
    "rendimento-hora médio real habitual de todos os trabalhos(r$/hora)",
    "índice gini",
    "população desocupada",
    "nível de ocupação",
    "população subutilizada",
    "taxa total mortalidade",
    "total pessoas ocupadas(1 000 pessoas)",
    "70 anos ou mais",
    "domicílio alugado",
    "domicílio cedido de outra forma",
    "taxa de desocupação",
    "45 a 59 anos",
    "taxa composta de subutilização",
    "população na força de trabalho potencial",
    "domicílio próprio - já pago",
    "60 a 69 anos",
    "60 anos ou mais",
    "domicílio cedido por familiar",
    "população ocupada",
    "15 a 29 anos",
    "domicílio próprio - pagando",
    "taxa de participação",
    "número mensal médio de leitos de internação (total)",
    "0 a 14 anos",
    "população ocupada em trabalhos formais",
    "saneamento basico total(1 000 pessoas)",
    "população",
    "população na força de trabalho",
    "total pessoas por condição de ocupação a domicílio(1 000 pessoas)",
    "população em idade de trabalhar",
    "30 a 44 anos"
)

for ($i = 0; $i -lt $labels.Length; $i++) {
    $row = $i + 2
    $value = $ws.Cells.Item($row, 1).Value2
    $ws.Cells.Item($row, 2).Value = $value
    $ws.Cells.Item($row, 1).Value = $labels[$i]
}
